# Updates the cryptocurrency snapshot table on Sheet1 (columns B-E, rows 2-51)
# to the refreshed values from the Thu Oct 24 16:12:00 UTC 2024 GitHub Actions run.
#
# Every write goes through a small helper that forces the destination cell to
# plain text first (matching the workbooks existing inline-string cells) and
# then restores the "Normal" style so numeric-looking text (e.g. "1.00", "8.00")
# is not silently reinterpreted as a number and does not pick up a stray format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value2 = $text
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue 2 4 "67.571.32"

# Row 3
Set-TextValue 3 4 "2.522.74"
Set-TextValue 3 5 "  -0.41%  "

# Row 4
Set-TextValue 4 5 "  -0.01%  "

# Row 5
Set-TextValue 5 4 "591.99"
Set-TextValue 5 5 "  +2.37%  "

# Row 6
Set-TextValue 6 4 "175.35"
Set-TextValue 6 5 "  +4.81%  "

# Row 7
Set-TextValue 7 5 "  +0.00%  "

# Row 8
Set-TextValue 8 5 "  +1.05%  "

# Row 9
Set-TextValue 9 4 "2.522.70"
Set-TextValue 9 5 "  -0.37%  "

# Row 10
Set-TextValue 10 4 "0.140"
Set-TextValue 10 5 "  +1.26%  "

# Row 11
Set-TextValue 11 5 "  +2.44%  "

# Row 12
Set-TextValue 12 5 "  +0.68%  "

# Row 13
Set-TextValue 13 4 "0.343"
Set-TextValue 13 5 "  -1.79%  "

# Row 14
Set-TextValue 14 5 "  +1.49%  "

# Row 15
Set-TextValue 15 4 "2.985.11"
Set-TextValue 15 5 "  -0.46%  "

# Row 16
Set-TextValue 16 5 "  +1.34%  "

# Row 17
Set-TextValue 17 4 "67.435.81"
Set-TextValue 17 5 "  +2.43%  "

# Row 18
Set-TextValue 18 4 "2.523.95"
Set-TextValue 18 5 "  -1.45%  "

# Row 19
Set-TextValue 19 4 "8.00"
Set-TextValue 19 5 "  +5.27%  "

# Row 20
Set-TextValue 20 5 "  +1.33%  "

# Row 21
Set-TextValue 21 4 "359.55"
Set-TextValue 21 5 "  +4.00%  "

# Row 22
Set-TextValue 22 4 "4.19"
Set-TextValue 22 5 "  +0.26%  "

# Row 23
Set-TextValue 23 4 "4.63"

# Row 24
Set-TextValue 24 4 "1.98"
Set-TextValue 24 5 "  +3.05%  "

# Row 25
Set-TextValue 25 5 "  +0.11%  "

# Row 26
Set-TextValue 26 4 "70.96"
Set-TextValue 26 5 "  +2.88%  "

# Row 27
Set-TextValue 27 4 "10.23"
Set-TextValue 27 5 "  +2.95%  "

# Row 28
Set-TextValue 28 4 "0.997"
Set-TextValue 28 5 "  -0.49%  "

# Row 29
Set-TextValue 29 4 "2.656.93"
Set-TextValue 29 5 "  -1.02%  "

# Row 30
Set-TextValue 30 5 "  +1.32%  "

# Row 31
Set-TextValue 31 4 "550.21"
Set-TextValue 31 5 "  +5.22%  "

# Row 32
Set-TextValue 32 4 "8.26"
Set-TextValue 32 5 "  +1.65%  "

# Row 33
Set-TextValue 33 4 "1.34"
Set-TextValue 33 5 "  +2.82%  "

# Row 34
Set-TextValue 34 4 "1.87"
Set-TextValue 34 5 "  +3.11%  "

# Row 35
Set-TextValue 35 4 "0.130"
Set-TextValue 35 5 "  -0.11%  "

# Row 36
Set-TextValue 36 5 "  +0.07%  "

# Row 37
Set-TextValue 37 5 "  +1.82%  "

# Row 38
Set-TextValue 38 4 "155.71"
Set-TextValue 38 5 "  -0.97%  "

# Row 39
Set-TextValue 39 4 "18.73"
Set-TextValue 39 5 "  +0.46%  "

# Row 40
Set-TextValue 40 4 "18.61"
Set-TextValue 40 5 "  +1.82%  "

# Row 41
Set-TextValue 41 4 "0.355"
Set-TextValue 41 5 "  +0.40%  "

# Row 42
Set-TextValue 42 5 "  +3.12%  "

# Row 43
Set-TextValue 43 4 "5.18"
Set-TextValue 43 5 "  +2.33%  "

# Row 44
Set-TextValue 44 2 "dogwifhat"
Set-TextValue 44 3 "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue 44 4 "2.52"
Set-TextValue 44 5 "  +5.50%  "

# Row 45
Set-TextValue 45 2 "USDe"
Set-TextValue 45 3 "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue 45 4 "1.00"
Set-TextValue 45 5 "  +0.01%  "

# Row 46
Set-TextValue 46 2 "ARBITRUM"
Set-TextValue 46 3 "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue 46 4 "0.562"
Set-TextValue 46 5 "  +1.19%  "

# Row 47
Set-TextValue 47 2 "BabyDogeCoin"
Set-TextValue 47 3 "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue 47 4 "0.0₆0281"
Set-TextValue 47 5 "  +0.20%  "

# Row 48
Set-TextValue 48 2 "Aave"
Set-TextValue 48 3 "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue 48 4 "146.86"
Set-TextValue 48 5 "  -0.20%  "

# Row 49
Set-TextValue 49 5 "  +1.24%  "

# Row 50
Set-TextValue 50 5 "  -1.05%  "

# Row 51
Set-TextValue 51 4 "0.0757"
Set-TextValue 51 5 "  +0.27%  "

